# BB_Build.xlsx edit: add Icon / ShowName columns (G/H), push Desc to I.
# Table becomes: ID, Type, SubType, Prefab, NormalStateFunc, UpStateFunc,
#                Icon, ShowName, Desc

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
# G1 used to hold "Desc"; it now becomes "Icon", a new "ShowName" header
# goes in H1, and "Desc" moves out to I1.
$ws.Cells.Item(1, 9).Value2 = $ws.Cells.Item(1, 7).Value2   # I1 = old G1 ("Desc")
$ws.Cells.Item(1, 7).Value2 = "Icon"                        # G1 = "Icon"
$ws.Cells.Item(1, 8).Value2 = "ShowName"                    # H1 = "ShowName"

# --- Data rows (2-10): short prefab names for Icon/ShowName ------------
$shortNames = @{
    2  = "Altar_1_1"
    3  = "Arena_1_1"
    4  = "Camp_1_1"
    5  = "GoldMine_1_1"
    6  = "Item_hourse_1_1"
    7  = "League_1_1"
    8  = "MagicHourse_1_1"
    9  = "Tower_1_1"
    10 = "Town_1_1"
}

for ($r = 2; $r -le 10; $r++) {
    # Move the existing Desc text (old column G) out to I, keep it as text.
    $descValue = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 9).NumberFormat = "@"
    $ws.Cells.Item($r, 9).Value2 = $descValue

    # New G = short prefab name (Icon)
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value2 = $shortNames[$r]

    # New H = ShowName, same text as Desc
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value2 = $descValue
}

# --- Column widths: columns G..I all width 11 (matches ColumnWidth=10.29 char units) ---
$ws.Columns.Item(7).ColumnWidth = 10.29
$ws.Columns.Item(8).ColumnWidth = 10.29
$ws.Columns.Item(9).ColumnWidth = 10.29

# --- Selection moves to H10 ---------------------------------------------
[void]$ws.Range("H10").Select()
